$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-style pattern: some "Price" values in column D look like plain numbers
# (e.g. "603.94"). Excel's normal Value-assignment auto-converts such strings
# to numeric cells, which would change both the stored type and the text
# representation (trailing float noise). Force these specific cells to stay
# text by flipping NumberFormat to Text ("@") for the assignment, then
# restoring the cell style to Normal so no stray formatting remains.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# --- Simple value updates (single cell changes) ---
Set-TextValue "D2" "65.346.47"
Set-TextValue "D3" "2.646.10"
$ws.Range("E4").Value  = "  -0.01%  "
Set-TextValue "D5" "603.94"
$ws.Range("E5").Value  = "  +2.72%  "
Set-TextValue "D6" "156.97"
$ws.Range("E6").Value  = "  +5.33%  "
$ws.Range("E7").Value  = "  -0.03%  "
$ws.Range("E8").Value  = "  +1.10%  "
$ws.Range("E9").Value  = "  +10.91%  "
$ws.Range("E10").Value = "  +6.18%  "
Set-TextValue "D11" "5.81"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("E12").Value = "  +2.70%  "
Set-TextValue "D13" "29.39"
$ws.Range("E13").Value = "  +6.79%  "
Set-TextValue "D14" "0.0000190"
$ws.Range("E14").Value = "  +22.66%  "
Set-TextValue "D15" "3.122.32"
$ws.Range("E15").Value = "  +1.90%  "
Set-TextValue "D16" "65.183.06"
$ws.Range("E16").Value = "  +3.09%  "
Set-TextValue "D17" "2.642.00"
$ws.Range("E17").Value = "  +1.57%  "
Set-TextValue "D18" "12.71"
$ws.Range("E18").Value = "  +5.52%  "
$ws.Range("E19").Value = "  +5.10%  "
Set-TextValue "D20" "360.20"
$ws.Range("E20").Value = "  +4.76%  "
$ws.Range("E21").Value = "  +8.19%  "
$ws.Range("E22").Value = "  -0.04%  "
Set-TextValue "D23" "69.27"
$ws.Range("E23").Value = "  +4.04%  "
$ws.Range("E24").Value = "  -0.53%  "
Set-TextValue "D25" "9.46"
$ws.Range("E25").Value = "  +3.39%  "
Set-TextValue "D27" "8.31"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("E29").Value = "  +15.07%  "

# --- Row 30/31 swap: Bittensor <-> PancakeSwap (with updated values) ---
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "2.20"
$ws.Range("E30").Value = "  +8.76%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D31" "544.35"
$ws.Range("E31").Value = "  -1.41%  "

$ws.Range("E32").Value = "  +0.01%  "
Set-TextValue "D33" "1.81"
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("E34").Value = "  +6.14%  "
$ws.Range("E35").Value = "  +6.29%  "
$ws.Range("E36").Value = "  +4.71%  "
Set-TextValue "D37" "20.63"
$ws.Range("E37").Value = "  +6.61%  "
$ws.Range("E38").Value = "  +5.32%  "
Set-TextValue "D39" "162.41"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  -0.02%  "
Set-TextValue "D41" "0.999"
$ws.Range("E41").Value = "  +0.04%  "
Set-TextValue "D42" "42.73"
$ws.Range("E42").Value = "  +8.19%  "
Set-TextValue "D43" "166.76"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("E44").Value = "  +4.69%  "
Set-TextValue "D45" "0.0618"
$ws.Range("E45").Value = "  +7.25%  "

# --- Row 46/47 swap: InjectiveProtocol <-> dogwifhat (with updated values) ---
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D46" "2.29"
$ws.Range("E46").Value = "  +9.85%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D47" "23.29"
$ws.Range("E47").Value = "  +1.79%  "

# --- Row 48/49 swap: VeChain <-> Mantle (with updated values) ---
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D48" "0.657"
$ws.Range("E48").Value = "  +4.06%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D49" "0.0264"
$ws.Range("E49").Value = "  +7.13%  "

$ws.Range("E50").Value = "  +2.88%  "
Set-TextValue "D51" "19.75"
$ws.Range("E51").Value = "  +3.94%  "
